# Workbook: ESP/data/service.xlsx
# "Added plos one to revs, youtube link for brazilian talk, and reorganised folders"
#
# This sheet (Hoja1) lists review/editorial activities. The relevant part of
# the edit re-orders the list of journals in column E (rows 6-16) and inserts
# a brand-new entry for "PLOS ONE", which pushes the closing rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-point the existing journal entries one slot further down the list ---
$ws.Range("E6").Value  = "Royal Society Open Science"
$ws.Range("E13").Value = "Basic and Applied Social Psychology"
$ws.Range("E14").Value = "Frontiers in Psychology"
$ws.Range("E15").Value = "Evolutionary Psychology"
$ws.Range("E16").Value = "Human Ethology Bulletin"

# --- New journal added to the reviewer list: PLOS ONE (row 7) ---
$ws.Range("E7").Value = "PLOS ONE"

# --- Insert a brand-new row (17) for "Summa Psicológica", pushing the final
#     two rows (Evaluación de Publicaciones..., Universidad Nacional de
#     Colombia) down from 17-18 to 18-19 ---
$ws.Rows.Item(17).Insert()
$ws.Range("E17").Value = "Summa Psicológica"
$ws.Range("E17").WrapText = $true

# --- Update the view: drop the frozen/scrolled topLeftCell and move the
#     active selection to reflect the newly-added row ---
$ws.Range("E24").Select() | Out-Null
